# Restored from revision of admin on 10/22/2020 08:26:38 AM.TEST Author: admin. Type: SAVE.
# Change the "From" threshold of the R30 ("Good Evening") rule row from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
